# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# to match the freshly scraped data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 481
$ws1.Range("F3").Value = 166
$ws1.Range("F4").Value = 8624
$ws1.Range("F5").Value = 106
$ws1.Range("F11").Value = 177
$ws1.Range("F13").Value = 466
$ws1.Range("F14").Value = 74
$ws1.Range("F17").Value = 6084
$ws1.Range("F18").Value = 201
$ws1.Range("F19").Value = 300
$ws1.Range("F20").Value = 2198
$ws1.Range("F21").Value = 80
$ws1.Range("F22").Value = 135
$ws1.Range("F24").Value = 442

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 481
$ws4.Range("F3").Value = 166
$ws4.Range("F4").Value = 8624
$ws4.Range("F5").Value = 106
$ws4.Range("F13").Value = 177
$ws4.Range("F15").Value = 466
$ws4.Range("F16").Value = 74
$ws4.Range("F20").Value = 6084
$ws4.Range("F22").Value = 201
$ws4.Range("F23").Value = 300
$ws4.Range("F24").Value = 2199
$ws4.Range("F25").Value = 80
$ws4.Range("F26").Value = 135
$ws4.Range("F28").Value = 442
